# Add two new condition rows to the "condition" sheet: itemCanEquipToRole
# and itemAlreadyEquiped — mirrors the source-table edit for the equip
# window cleanup (remove listener, finish window related).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A26").Value = "itemCanEquipToRole"
$ws.Range("B26").Value = "某角色可以装备"
$ws.Range("C26").Value = "data"
$ws.Range("D26").Value = "role"
$ws.Range("E26").Value = "canEquip:"
$ws.Range("F26").Value = "data"
$ws.Range("G26").Value = "item"

$ws.Range("A27").Value = "itemAlreadyEquiped"
$ws.Range("B27").Value = "已经装备上了"
$ws.Range("C27").Value = "data"
$ws.Range("D27").Value = "item"
$ws.Range("E27").Value = "isEquiped"
$ws.Range("F27").Value = ";"
$ws.Range("G27").Value = ";"

$ws.Range("A27").Select() | Out-Null
